$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds the "last changed" date stamp for every data
# row (rows 2-242). The automated update bumps this date by one day.
$ws.Range("C2:C242").Value = 46061
